# Updated cryptos list on Fri Jul 19 17:29:35 UTC 2024 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures scraped from
# coinranking.com; two coin pairs (InjectiveProtocol/EnergySwap and OKB/Maker)
# swapped rank positions, so their Coin/Link/Price/Volume cells are rewritten too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.569.01'
$ws.Range("E2").Value = '  +4.31%  '
$ws.Range("D3").Value = '3.505.69'
$ws.Range("E3").Value = '  +2.61%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.71'
$ws.Range("E5").Value = '  +3.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.76'
$ws.Range("E6").Value = '  +7.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.510.06'
$ws.Range("E8").Value = '  +2.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.28'
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +4.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.437'
$ws.Range("E12").Value = '  +2.06%  '
$ws.Range("D13").Value = '4.109.13'
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.135'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.09'
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").Value = '66.569.84'
$ws.Range("E16").Value = '  +4.18%  '
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '3.501.67'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.29'
$ws.Range("E19").Value = '  +2.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.06'
$ws.Range("E20").Value = '  +2.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.82'
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.00'
$ws.Range("E22").Value = '  +3.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.95'
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.529'
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("E26").Value = '  +6.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("E27").Value = '  +4.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.181'
$ws.Range("E28").Value = '  +1.87%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.40'
$ws.Range("E30").Value = '  +4.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.48'
$ws.Range("E31").Value = '  +5.82%  '
$ws.Range("E32").Value = '  +3.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.48'
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.41'
$ws.Range("E34").Value = '  +5.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.55'
$ws.Range("E36").Value = '  +1.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.65'
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.904'
$ws.Range("E38").Value = '  +8.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  +5.27%  '
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.48'
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.23'
$ws.Range("E42").Value = '  +5.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.69'
$ws.Range("E43").Value = '  +4.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.61'
$ws.Range("E44").Value = '  +4.39%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.808.47'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.53'
$ws.Range("E46").Value = '  +1.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0315'
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '357.18'
$ws.Range("E48").Value = '  +8.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.52'
$ws.Range("E49").Value = '  +8.40%  '
$ws.Range("E50").Value = '  +5.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '32.73'
$ws.Range("E51").Value = '  +8.72%  '
